# Applies the diff: updates the date header and the 25 division-problem
# answer cells in the single table, addressed by (row, col) to avoid
# ambiguity since several old/new answer strings repeat.

$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-09-27 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-28 Thursday", 2) | Out-Null

$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "80÷3=26, 2"
$t.Cell(1,2).Range.Text = "72÷3=24, 0"
$t.Cell(1,3).Range.Text = "47÷8=5, 7"
$t.Cell(1,4).Range.Text = "49÷3=16, 1"
$t.Cell(1,5).Range.Text = "13÷6=2, 1"
$t.Cell(5,1).Range.Text = "49÷4=12, 1"
$t.Cell(5,2).Range.Text = "56÷5=11, 1"
$t.Cell(5,3).Range.Text = "88÷6=14, 4"
$t.Cell(5,4).Range.Text = "49÷4=12, 1"
$t.Cell(5,5).Range.Text = "57÷4=14, 1"
$t.Cell(9,1).Range.Text = "37÷5=7, 2"
$t.Cell(9,2).Range.Text = "47÷9=5, 2"
$t.Cell(9,3).Range.Text = "13÷4=3, 1"
$t.Cell(9,4).Range.Text = "64÷2=32, 0"
$t.Cell(9,5).Range.Text = "94÷9=10, 4"
$t.Cell(13,1).Range.Text = "89÷9=9, 8"
$t.Cell(13,2).Range.Text = "41÷3=13, 2"
$t.Cell(13,3).Range.Text = "81÷5=16, 1"
$t.Cell(13,4).Range.Text = "78÷2=39, 0"
$t.Cell(13,5).Range.Text = "18÷2=9, 0"
$t.Cell(17,1).Range.Text = "47÷2=23, 1"
$t.Cell(17,2).Range.Text = "66÷9=7, 3"
$t.Cell(17,3).Range.Text = "14÷7=2, 0"
$t.Cell(17,4).Range.Text = "17÷3=5, 2"
$t.Cell(17,5).Range.Text = "39÷3=13, 0"
